# Apply the Saldo.xlsx update.
# Work from the bottom of the sheet upward so that earlier row numbers
# are not invalidated by later deletes/inserts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33: 004855960 CLERIA 792.77  -> removed (this account/name is
# consolidated into the new row 2 below with an updated balance).
$ws.Rows.Item(33).Delete()

# Row 16: 004497825 PRISCILLA 3084.67 -> removed.
$ws.Rows.Item(16).Delete()

# Insert a new row at position 15 for 005061124 / BRUNO / 4702.02
$ws.Rows.Item(15).Insert()
$ws.Cells.Item(15, 1).NumberFormat = "@"
$ws.Cells.Item(15, 1).Value = "005061124"
$ws.Cells.Item(15, 2).Value = "BRUNO"
$ws.Cells.Item(15, 3).Value = 4702.02

# Row 13: 004646727 RENATA 7683.85 -> 004472386 GABRIEL 13355.32
$ws.Cells.Item(13, 1).NumberFormat = "@"
$ws.Cells.Item(13, 1).Value = "004472386"
$ws.Cells.Item(13, 2).Value = "GABRIEL"
$ws.Cells.Item(13, 3).Value = 13355.32

# Row 8: 005061124 BRUNO 23724.23 (old, unrelated record) -> removed.
$ws.Rows.Item(8).Delete()

# Row 7: 004474776 GILSON 23887.9 -> removed.
$ws.Rows.Item(7).Delete()

# Row 6: 004472386 GABRIEL 33620.62 -> 005366671 TATIANA 35000
$ws.Cells.Item(6, 1).NumberFormat = "@"
$ws.Cells.Item(6, 1).Value = "005366671"
$ws.Cells.Item(6, 2).Value = "TATIANA"
$ws.Cells.Item(6, 3).Value = 35000

# Row 5: 004213929 RODOLFO 50000 -> 004213929 RODOLFO 44542
$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = "004213929"
$ws.Cells.Item(5, 2).Value = "RODOLFO"
$ws.Cells.Item(5, 3).Value = 44542

# Row 3: 005381719 MARIA 73733.2 -> 004804125 EDUARDO 93661.99
$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "004804125"
$ws.Cells.Item(3, 2).Value = "EDUARDO"
$ws.Cells.Item(3, 3).Value = 93661.99

# Row 2: 003987275 ELIANE 113108.96 -> 004855960 CLERIA 209277.61
$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "004855960"
$ws.Cells.Item(2, 2).Value = "CLERIA"
$ws.Cells.Item(2, 3).Value = 209277.61
